# feat: Add keyboard shortcuts
# New catalog entry "Iditarod" / "Soon Hee Newbold" / grade "4" is inserted
# as the new first row of the table; the existing "Sakura" row shifts down
# to row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record by inserting a row above the current row 1,
# which pushes the existing data down to row 2.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Iditarod"
$ws.Range("B1").Value = "Soon Hee Newbold"
$ws.Range("C1").Value = "'"
$ws.Range("D1").Value = "'"
# Force this cell to be stored as text (matching the rest of the column,
# e.g. row 2's "5") rather than being auto-converted to a number.
$ws.Range("E1").Value = "'4"
$ws.Range("F1").Value = "'"
